$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.190.67'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.859.54'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7092'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.29'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3073'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07641'
$ws.Range('E9').Value = '  -3.23%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '24.72'
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08384'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.860.81'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.179'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7074'
$ws.Range('E14').Value = '  -3.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.17'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.171.52'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.926'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '242.77'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007808'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.113.63'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.11'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.846'
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1591'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.59'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.877'
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.48'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.320'
$ws.Range('E29').Value = '  -3.09%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.387'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.225'
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05114'
$ws.Range('E33').Value = '  -3.74%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.8052'
$ws.Range('E34').Value = '  +10.96%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.916'
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.164'
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.683'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01840'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.686'
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.178.78'
$ws.Range('E40').Value = '  -7.01%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8922'
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '72.56'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '101.68'
$ws.Range('E45').Value = '  -1.93%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.012.07'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5155'
$ws.Range('E47').Value = '  -3.24%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.776'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.234'
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.9966'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4255'
$ws.Range('E51').Value = '  -1.82%  '
